$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 321, shifting the existing rows 321-367 down to 322-368.
$ws.Rows.Item(321).EntireRow.Insert()

# Populate the new row 321 with the new "Ajo" price record.
$ws.Range("A321").Value = 6
$ws.Range("B321").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C321").Value = "Metropolitana"
$ws.Range("D321").Value = 44491
$ws.Range("E321").Value = 13
$ws.Range("F321").Value = 100112003
$ws.Range("G321").Value = "Ajo"
$ws.Range("H321").Value = "Chino"
$ws.Range("I321").Value = "Primera"
$ws.Range("J321").Value = 2200
$ws.Range("K321").Value = 16500
$ws.Range("L321").Value = 17000
$ws.Range("M321").Value = 16705
$ws.Range("N321").Value = "$/caja 10 kilos"
$ws.Range("O321").Value = "China"
$ws.Range("P321").Value = 1670
$ws.Range("Q321").Value = 10
$ws.Range("R321").Value = "Hortaliza"
